$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: "time_taken" — copy formatting from the neighboring
# header cell (E1: bold, bordered, centered) so it matches the other headers.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# New data column values (F2:F5), left with default (unstyled) formatting,
# consistent with the other data cells in the sheet.
$ws.Range("F2").Value = "2021-10-05 13:38:48.380330"
$ws.Range("F3").Value = "2021-10-05 13:38:48.380337"
$ws.Range("F4").Value = "2021-10-05 13:38:48.380340"
$ws.Range("F5").Value = "2021-10-05 13:38:48.380342"

$excel.CutCopyMode = $false
